$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.548.55'
$ws.Range("E2").Value = '  -1.44%  '
$ws.Range("D3").Value = '1.666.37'
$ws.Range("E3").Value = '  -3.40%  '
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '215.25'
$ws.Range("E5").Value = '  -1.54%  '
$ws.Range("E6").Value = '  -1.78%  '
$ws.Range("E7").Value = '  +0.06%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '23.54'
$ws.Range("E8").Value = '  -1.73%  '
$ws.Range("E9").Value = '  -1.15%  '
$ws.Range("E10").Value = '  -1.81%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0879'
$ws.Range("E11").Value = '  -2.42%  '
$ws.Range("D12").Value = '1.901.86'
$ws.Range("E12").Value = '  -3.44%  '
$ws.Range("D13").Value = '1.666.19'
$ws.Range("E13").Value = '  -3.62%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.15'
$ws.Range("E14").Value = '  -2.47%  '
$ws.Range("E15").Value = '  -2.06%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '66.27'
$ws.Range("E16").Value = '  -2.42%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '249.19'
$ws.Range("E17").Value = '  +2.35%  '
$ws.Range("D18").Value = '27.568.89'
$ws.Range("E18").Value = '  -1.19%  '
$ws.Range("D19").Value = '0.0₃0733'
$ws.Range("E19").Value = '  -2.93%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.54'
$ws.Range("E20").Value = '  -4.59%  '
$ws.Range("E21").Value = '  -0.03%  '
$ws.Range("E22").Value = '  -3.08%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.29'
$ws.Range("E23").Value = '  -4.52%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.03'
$ws.Range("E24").Value = '  -5.61%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '146.55'
$ws.Range("E25").Value = '  -1.82%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '16.59'
$ws.Range("E26").Value = '  -1.28%  '
$ws.Range("E27").Value = '  -4.86%  '
$ws.Range("E28").Value = '  -2.20%  '
$ws.Range("E29").Value = '  +0.15%  '
$ws.Range("E30").Value = '  +3.82%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0510'
$ws.Range("E32").Value = '  -2.53%  '
$ws.Range("D33").Value = '1.474.54'
$ws.Range("E33").Value = '  -0.76%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.11'
$ws.Range("E34").Value = '  -5.38%  '
$ws.Range("E35").Value = '  -5.38%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.940'
$ws.Range("E36").Value = '  -1.95%  '
$ws.Range("E37").Value = '  -1.05%  '
$ws.Range("E38").Value = '  -5.83%  '
$ws.Range("E39").Value = '  -2.22%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '69.75'
$ws.Range("E40").Value = '  -2.16%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.00'
$ws.Range("E41").Value = '  +0.04%  '
$ws.Range("E42").Value = '  -4.93%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.43'
$ws.Range("E43").Value = '  -6.96%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.21'
$ws.Range("E44").Value = '  -3.23%  '
$ws.Range("D45").Value = '1.810.14'
$ws.Range("E45").Value = '  -3.40%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.791'
$ws.Range("E46").Value = '  -0.25%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.71'
$ws.Range("E47").Value = '  -1.61%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '89.54'
$ws.Range("E48").Value = '  -1.82%  '
$ws.Range("E49").Value = '  -1.30%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '42.04'
$ws.Range("E50").Value = '  +16.20%  '
$ws.Range("E51").Value = '  -3.17%  '
